# "fixed error in experiment format"
#
# Survey 2 was missing the "Pseudo-Random Question Width:" label row that
# Survey 1 already has right after the "Questions Per Page:" row (row 7).
# Insert the matching row into Survey 2, bring it into focus, and select it.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Survey 2")

# Shift rows 8..68 down to 9..69 and create a fresh row 8.
$ws2.Rows.Item(8).Insert()

# Populate the new row the same way the equivalent row on "Survey 1" (row 8)
# is populated: label in column A, empty formatted cell in column B.
$ws2.Cells.Item(8, 1).Value = "Pseudo-Random Question Width:"

# Make "Survey 2" the active sheet and select the newly inserted row.
$ws2.Activate()
$ws2.Range("A8:D8").Select()
